# Rotate the fact rows: the observation that was in row 5 moves up to row 2,
# and the former rows 2-4 shift down to rows 3-5 (row 5's old data becomes
# the new row 2; everything else shifts down by one row).
#
# Columns C, K, P, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY are identical
# across rows 2-5 already, so only the columns below actually change value
# as a result of the rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New row 2 (was row 5) ----
$ws.Range("A2").Value = 111645826
$ws.Range("B2").Value = 94134
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 53
$ws.Range("F2").Value = "Vedtrappmossa"
$ws.Range("G2").Value = "Crossocalyx hellerianus"
$ws.Range("H2").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("I2").ClearContents()
$ws.Range("Q2").Value = 369470
$ws.Range("R2").Value = 6635346
$ws.Range("S2").Value = 10
$ws.Range("Z2").Value = "13:25"
$ws.Range("AB2").Value = "13:25"

# ---- New row 3 (was row 2) ----
$ws.Range("A3").Value = 111644287
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
# "35" is stored as text (not a number) in the source data, so force text
# formatting before assigning, then restore the default style so no new
# cell format is left behind.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "35"
$ws.Range("I3").Style = "Normal"
$ws.Range("Q3").Value = 369410
$ws.Range("R3").Value = 6635288
$ws.Range("S3").Value = 10
$ws.Range("Z3").Value = "12:37"
$ws.Range("AB3").Value = "12:37"

# ---- New row 4 (was row 3) ----
# (I4 was already empty and stays empty - no change needed there.)
$ws.Range("A4").Value = 111644956
$ws.Range("B4").Value = 96348
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("Q4").Value = 369440
$ws.Range("R4").Value = 6635308
$ws.Range("S4").Value = 10
$ws.Range("Z4").Value = "12:52"
$ws.Range("AB4").Value = "12:52"

# ---- New row 5 (was row 4) ----
# (I5 was already empty and stays empty - no change needed there.)
$ws.Range("A5").Value = 111644923
$ws.Range("B5").Value = 56543
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 103021
$ws.Range("F5").Value = "Talltita"
$ws.Range("G5").Value = "Poecile montanus"
$ws.Range("H5").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q5").Value = 369436
$ws.Range("R5").Value = 6635294
$ws.Range("S5").Value = 25
$ws.Range("Z5").Value = "12:52"
$ws.Range("AB5").Value = "12:52"
